$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value
$ws.Range("B11").Value = 5

# Update total correct count and corr/total marks text
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
